$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# E11 previously held the shared string "Good Night"; it now holds the
# literal text "= error(""fail"")" (entered as text, not as a formula -
# the leading apostrophe forces Excel to store it as a text value with a
# quote-prefix instead of evaluating it as "=error(...)").
$ws.Range("E11").Value = "'= error(""fail"")"

# Leave selection where the author left it when saving.
[void]$ws.Range("F9").Select()
